# Update the "WORCreation" sheet so the supplier/client referenced changes
# from "ABC News" to "HSBC".
#   A3 (worSupplier)  : "ABC News" -> "HSBC"
#   C3 (worsourceName): "Master Service Agreement - ABC News" -> "Master Service Agreement - HSBC"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("WORCreation")

$ws.Range("A3").Value = "HSBC"
$ws.Range("C3").Value = "Master Service Agreement - HSBC"
